# natmiOut/YoungD0/LR-pairs_lrc2p/Vcan-Itgb1.xlsx
# Rebuild the sending/target-cluster crosstab for the Vcan-Itgb1 ligand-receptor
# pair after adding the "ECs" sending cluster (per Dr Hou's advice): the result
# is a 3x3 (sending cluster x target cluster) grid, rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of rows; first element of each is the target worksheet row
# number, the rest are the column A..T values in order.
$data = @(
  , @(2,  "ECs", "Vcan", "Itgb1", "ECs", 2, 0.6666666666666666, 1.529781, 4.589343, 0.01315047351877542, 0.01315047351877542, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 172.120849327152, 1549.087643944368, 0.004307191657458678, 0.004307191657458678)
  , @(3,  "ECs", "Vcan", "Itgb1", "FAPs", 2, 0.6666666666666666, 1.529781, 4.589343, 0.01315047351877542, 0.01315047351877542, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 162.637850111946, 1463.740651007514, 0.004069886907528076, 0.004069886907528077)
  , @(4,  "ECs", "Vcan", "Itgb1", "sCs", 2, 0.6666666666666666, 1.529781, 4.589343, 0.01315047351877542, 0.01315047351877542, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 190.750925187482, 1716.758326687338, 0.004773394953788666, 0.004773394953788667)
  , @(5,  "FAPs", "Vcan", "Itgb1", "ECs", 3, 1, 103.676216, 311.028648, 0.8912330150752564, 0.8912330150752563, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 11664.96273188467, 104984.664586962, 0.29190670601353, 0.29190670601353)
  , @(6,  "FAPs", "Vcan", "Itgb1", "FAPs", 3, 1, 103.676216, 311.028648, 0.8912330150752564, 0.8912330150752563, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 11022.28154094066, 99200.53386846589, 0.2758241043132664, 0.2758241043132664)
  , @(7,  "FAPs", "Vcan", "Itgb1", "sCs", 3, 1, 103.676216, 311.028648, 0.8912330150752564, 0.8912330150752563, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 12927.55899173622, 116348.030925626, 0.32350220474846, 0.32350220474846)
  , @(8,  "sCs", "Vcan", "Itgb1", "ECs", 3, 1, 11.12297, 33.36891, 0.09561651140596822, 0.09561651140596822, 3, 1, 112.513392, 337.540176, 0.3275312977368564, 0.3275312977368564, 1251.48308381424, 11263.34775432816, 0.0313174000658677, 0.0313174000658677)
  , @(9,  "sCs", "Vcan", "Itgb1", "FAPs", 3, 1, 11.12297, 33.36891, 0.09561651140596822, 0.09561651140596822, 3, 1, 106.314466, 318.943398, 0.3094859589441663, 0.3094859589441664, 1182.53261588402, 10642.79354295618, 0.02959196772337189, 0.02959196772337189)
  , @(10, "sCs", "Vcan", "Itgb1", "sCs", 3, 1, 11.12297, 33.36891, 0.09561651140596822, 0.09561651140596822, 3, 1, 124.6916553333333, 374.074966, 0.3629827433189773, 0.3629827433189773, 1386.941541523007, 12482.47387370706, 0.03470714361672863, 0.03470714361672863)
)

foreach ($rowSpec in $data) {
  $r = $rowSpec[0]
  for ($i = 1; $i -lt $rowSpec.Length; $i++) {
    $ws.Cells.Item($r, $i).Value = $rowSpec[$i]
  }
}

Write-Host "Updated rows 2-10 for Vcan-Itgb1 ECs/FAPs/sCs crosstab"
